$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header and data (column E)
# Written in this particular order so the shared-strings table is built
# with the same index assignment as the source workbook.
$ws.Range("E1").Value = "TestCase_Label"
$ws.Range("E2").Value = "valid_email_valid_password"
$ws.Range("E5").Value = "valid_email_invalid_password"
$ws.Range("E3").Value = "invalid_email_invalid_password"
$ws.Range("E4").Value = "invalid_email_valid_password"
$ws.Range("E6").Value = "empty_email_empty_password"

# Set new column width (closest achievable value to the source width of
# 27.59765625 chars given this engine's column-width rounding granularity)
$ws.Columns.Item(5).ColumnWidth = 26.833333333333332

# Update selection to match the author's final cursor position
$ws.Range("E12").Select()
